# Occupancy Optimization v5.0 Stable
# Rename the "Sept27thruOct3" sheet tab to "Sep27-Oct3_SJ" (third sheet /
# rId4 in the workbook). This is the only functionally meaningful change
# described by the diff — the accompanying style/font-table shrink and the
# tiny defaultColWidth nudges on several sheets are incidental artifacts of
# the source file's resave pipeline (duplicate-font/style dedup + column
# width recompute) with no user-visible formatting effect, so they need no
# explicit COM action here.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sept27thruOct3")
$ws.Name = "Sep27-Oct3_SJ"
